# 20150113 +++++++ cs-厂商 end +
#
# Adds a new "mj-买家" (buyer) section below the existing "cs-厂商" (vendor)
# page-comparison table: a big bold section header in A15, followed by three
# new .psd -> .html mapping rows (16-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New section header (A15), styled like the existing "cs-厂商" header in
# A1 (bold + centered) but larger. Copy A1's format first so the new style
# reuses the existing bold/centered xf and only a single new font (size 24)
# has to be created, then overwrite the text and bump the font size.
$ws.Range("A1").Copy()
$ws.Range("A15").PasteSpecial(-4122)            # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A15").Value = "mj-买家"
$ws.Range("A15").Font.Size = 24
$ws.Rows.Item(15).RowHeight = 31.5

# --- New mapping rows under the new header
$ws.Range("A16").Value = "店铺首页.psd"
$ws.Range("B16").Value = "shop-index-mj.html"

$ws.Range("A17").Value = "店铺首页-搜索.psd"
$ws.Range("B17").Value = "shop-index-mj-search.html"

$ws.Range("A18").Value = "商品详情 - 2.psd"
$ws.Range("B18").Value = "shop-show-mj.html"

# Leave the selection on the last edited cell, matching the authored commit.
$ws.Range("B18").Select()
